# Typocrypha spellDictionary.xlsx update
# - Adds two new debug/utility spells ("death" and "stun") to the attack table,
#   pushing the later tables (element table, style table) down by two rows.
# - Moves the active cell selection to A10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows right before the current row 10 (the "END" marker
# that closes the attack-spell table), shifting everything below it down by
# two rows. Doing this twice, one row at a time, mirrors how Excel performs a
# two-row insert and keeps all existing formatting/styles on the rows below.
$ws.Rows.Item(10).Insert()
$ws.Rows.Item(10).Insert()

# New row 10: "death" debug spell (kills the caster).
$ws.Range("B10").Value = "attack"
$ws.Range("C10").Value = "kills self (debug)"
$ws.Range("D10").Value = 999
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 100
$ws.Range("G10").Value = 0
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = "mS"

# New row 11: "stun" debug spell (staggers all targets).
$ws.Range("A11").Value = "stun"
$ws.Range("B11").Value = "attack"
$ws.Range("C11").Value = "deal stagger to all"
$ws.Range("D11").Value = 1
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 100
$ws.Range("G11").Value = 100
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = "LMR"

# Set A10 last so the shared-string table gets the new entries appended in
# the same order as the reference edit (kills self (debug), mS, stun,
# deal stagger to all, LMR, death).
$ws.Range("A10").Value = "death"

# Update the saved selection/active cell to A10.
[void]$ws.Range("A10").Select()
